$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.1829625
$ws.Range("H2").Value = 68.365925
$ws.Range("I2").Value = 0.6188383653293237
$ws.Range("J2").Value = 0.5689764244710266
$ws.Range("M2").Value = 44.88644
$ws.Range("N2").Value = 89.77288
$ws.Range("O2").Value = 0.08941539400600468
$ws.Range("P2").Value = 0.06313518933231116
$ws.Range("Q2").Value = 1534.3514952785
$ws.Range("R2").Value = 6137.405981114
$ws.Range("S2").Value = 0.05533367626195334
$ws.Range("T2").Value = 0.0359224342845997
$ws.Range("G3").Value = 34.1829625
$ws.Range("H3").Value = 68.365925
$ws.Range("I3").Value = 0.6188383653293237
$ws.Range("J3").Value = 0.5689764244710266
$ws.Range("O3").Value = 0.06989553179102848
$ws.Range("P3").Value = 0.07402865606362163
$ws.Range("Q3").Value = 1199.394297917521
$ws.Range("R3").Value = 7196.365787505126
$ws.Range("S3").Value = 0.04325403663738384
$ws.Range("T3").Value = 0.04212056003547482
$ws.Range("G4").Value = 34.1829625
$ws.Range("H4").Value = 68.365925
$ws.Range("I4").Value = 0.6188383653293237
$ws.Range("J4").Value = 0.5689764244710266
$ws.Range("M4").Value = 29.15707
$ws.Range("N4").Value = 87.47121
$ws.Range("O4").Value = 0.05808192634815011
$ws.Range("P4").Value = 0.06151647807752574
$ws.Range("Q4").Value = 996.675030419875
$ws.Range("R4").Value = 5980.05018251925
$ws.Range("S4").Value = 0.03594332435646739
$ws.Range("T4").Value = 0.03500142574260089
$ws.Range("G5").Value = 34.1829625
$ws.Range("H5").Value = 68.365925
$ws.Range("I5").Value = 0.6188383653293237
$ws.Range("J5").Value = 0.5689764244710266
$ws.Range("M5").Value = 39.1954995
$ws.Range("N5").Value = 78.39099899999999
$ws.Range("O5").Value = 0.07807883697291786
$ws.Range("P5").Value = 0.05513057577983479
$ws.Range("Q5").Value = 1339.818289577269
$ws.Range("R5").Value = 5359.273158309074
$ws.Range("S5").Value = 0.04831817983913524
$ws.Range("T5").Value = 0.03136799788623938
$ws.Range("G6").Value = 34.1829625
$ws.Range("H6").Value = 68.365925
$ws.Range("I6").Value = 0.6188383653293237
$ws.Range("J6").Value = 0.5689764244710266
$ws.Range("M6").Value = 270.2169853333333
$ws.Range("N6").Value = 810.650956
$ws.Range("O6").Value = 0.5382819000726007
$ws.Range("P6").Value = 0.5701120604516535
$ws.Range("Q6").Value = 9236.817076512383
$ws.Range("R6").Value = 55420.9024590743
$ws.Range("S6").Value = 0.3331094911272905
$ws.Range("T6").Value = 0.3243803217035916
$ws.Range("G7").Value = 34.1829625
$ws.Range("H7").Value = 68.365925
$ws.Range("I7").Value = 0.6188383653293237
$ws.Range("J7").Value = 0.5689764244710266
$ws.Range("M7").Value = 83.455535
$ws.Range("N7").Value = 250.366605
$ws.Range("O7").Value = 0.1662464108092982
$ws.Range("P7").Value = 0.1760770402950531
$ws.Range("Q7").Value = 2852.757423322438
$ws.Range("R7").Value = 17116.54453993462
$ws.Range("S7").Value = 0.1028796571070933
$ws.Range("T7").Value = 0.1001836848185202
$ws.Range("I8").Value = 0.06738425137939692
$ws.Range("J8").Value = 0.09293230485581538
$ws.Range("M8").Value = 44.88644
$ws.Range("N8").Value = 89.77288
$ws.Range("O8").Value = 0.08941539400600468
$ws.Range("P8").Value = 0.06313518933231116
$ws.Range("Q8").Value = 167.0729105607067
$ws.Range("R8").Value = 1002.43746336424
$ws.Range("S8").Value = 0.00602518938688844
$ws.Range("T8").Value = 0.005867298662159964
$ws.Range("I9").Value = 0.06738425137939692
$ws.Range("J9").Value = 0.09293230485581538
$ws.Range("O9").Value = 0.06989553179102848
$ws.Range("P9").Value = 0.07402865606362163
$ws.Range("S9").Value = 0.004709858084503292
$ws.Range("T9").Value = 0.006879653633370792
$ws.Range("I10").Value = 0.06738425137939692
$ws.Range("J10").Value = 0.09293230485581538
$ws.Range("M10").Value = 29.15707
$ws.Range("N10").Value = 87.47121
$ws.Range("O10").Value = 0.05808192634815011
$ws.Range("P10").Value = 0.06151647807752574
$ws.Range("Q10").Value = 108.5262397357033
$ws.Range("R10").Value = 976.73615762133
$ws.Range("S10").Value = 0.003913807125643365
$ws.Range("T10").Value = 0.005716868094356706
$ws.Range("I11").Value = 0.06738425137939692
$ws.Range("J11").Value = 0.09293230485581538
$ws.Range("M11").Value = 39.1954995
$ws.Range("N11").Value = 78.39099899999999
$ws.Range("O11").Value = 0.07807883697291786
$ws.Range("P11").Value = 0.05513057577983479
$ws.Range("Q11").Value = 145.8905224461045
$ws.Range("R11").Value = 875.343134676627
$ws.Range("S11").Value = 0.005261283977994047
$ws.Range("T11").Value = 0.005123411475248239
$ws.Range("I12").Value = 0.06738425137939692
$ws.Range("J12").Value = 0.09293230485581538
$ws.Range("M12").Value = 270.2169853333333
$ws.Range("N12").Value = 810.650956
$ws.Range("O12").Value = 0.5382819000726007
$ws.Range("P12").Value = 0.5701120604516535
$ws.Range("Q12").Value = 1005.781216389176
$ws.Range("R12").Value = 9052.030947502588
$ws.Range("S12").Value = 0.03627172286747154
$ws.Range("T12").Value = 0.05298182780387011
$ws.Range("I13").Value = 0.06738425137939692
$ws.Range("J13").Value = 0.09293230485581538
$ws.Range("M13").Value = 83.455535
$ws.Range("N13").Value = 250.366605
$ws.Range("O13").Value = 0.1662464108092982
$ws.Range("P13").Value = 0.1760770402950531
$ws.Range("Q13").Value = 310.6318775748517
$ws.Range("R13").Value = 2795.686898173665
$ws.Range("S13").Value = 0.01120238993689624
$ws.Range("T13").Value = 0.01636324518680956
$ws.Range("G14").Value = 1.560491
$ws.Range("H14").Value = 4.681473
$ws.Range("I14").Value = 0.02825067311094296
$ws.Range("J14").Value = 0.03896162845449177
$ws.Range("M14").Value = 44.88644
$ws.Range("N14").Value = 89.77288
$ws.Range("O14").Value = 0.08941539400600468
$ws.Range("P14").Value = 0.06313518933231116
$ws.Range("Q14").Value = 70.04488564204
$ws.Range("R14").Value = 420.26931385224
$ws.Range("S14").Value = 0.002526045067149807
$ws.Range("T14").Value = 0.0024598497891695
$ws.Range("G15").Value = 1.560491
$ws.Range("H15").Value = 4.681473
$ws.Range("I15").Value = 0.02825067311094296
$ws.Range("J15").Value = 0.03896162845449177
$ws.Range("O15").Value = 0.06989553179102848
$ws.Range("P15").Value = 0.07402865606362163
$ws.Range("Q15").Value = 54.75370975677167
$ws.Range("R15").Value = 492.7833878109451
$ws.Range("S15").Value = 0.001974595820543867
$ws.Range("T15").Value = 0.002884276992536185
$ws.Range("G16").Value = 1.560491
$ws.Range("H16").Value = 4.681473
$ws.Range("I16").Value = 0.02825067311094296
$ws.Range("J16").Value = 0.03896162845449177
$ws.Range("M16").Value = 29.15707
$ws.Range("N16").Value = 87.47121
$ws.Range("O16").Value = 0.05808192634815011
$ws.Range("P16").Value = 0.06151647807752574
$ws.Range("Q16").Value = 45.49934532137
$ws.Range("R16").Value = 409.49410789233
$ws.Range("S16").Value = 0.001640853514915454
$ws.Range("T16").Value = 0.002396782162685446
$ws.Range("G17").Value = 1.560491
$ws.Range("H17").Value = 4.681473
$ws.Range("I17").Value = 0.02825067311094296
$ws.Range("J17").Value = 0.03896162845449177
$ws.Range("M17").Value = 39.1954995
$ws.Range("N17").Value = 78.39099899999999
$ws.Range("O17").Value = 0.07807883697291786
$ws.Range("P17").Value = 0.05513057577983479
$ws.Range("Q17").Value = 61.1642242102545
$ws.Range("R17").Value = 366.985345261527
$ws.Range("S17").Value = 0.00220577970020451
$ws.Range("T17").Value = 0.002147977010016126
$ws.Range("G18").Value = 1.560491
$ws.Range("H18").Value = 4.681473
$ws.Range("I18").Value = 0.02825067311094296
$ws.Range("J18").Value = 0.03896162845449177
$ws.Range("M18").Value = 270.2169853333333
$ws.Range("N18").Value = 810.650956
$ws.Range("O18").Value = 0.5382819000726007
$ws.Range("P18").Value = 0.5701120604516535
$ws.Range("Q18").Value = 421.6711736597986
$ws.Range("R18").Value = 3795.040562938188
$ws.Range("S18").Value = 0.01520682600048831
$ws.Range("T18").Value = 0.02221249427674208
$ws.Range("G19").Value = 1.560491
$ws.Range("H19").Value = 4.681473
$ws.Range("I19").Value = 0.02825067311094296
$ws.Range("J19").Value = 0.03896162845449177
$ws.Range("M19").Value = 83.455535
$ws.Range("N19").Value = 250.366605
$ws.Range("O19").Value = 0.1662464108092982
$ws.Range("P19").Value = 0.1760770402950531
$ws.Range("Q19").Value = 130.231611267685
$ws.Range("R19").Value = 1172.084501409165
$ws.Range("S19").Value = 0.004696573007641017
$ws.Range("T19").Value = 0.006860248223342433
$ws.Range("G20").Value = 11.3729585
$ws.Range("H20").Value = 22.745917
$ws.Range("I20").Value = 0.2058927176688748
$ws.Range("J20").Value = 0.1893032314851988
$ws.Range("M20").Value = 44.88644
$ws.Range("N20").Value = 89.77288
$ws.Range("O20").Value = 0.08941539400600468
$ws.Range("P20").Value = 0.06313518933231116
$ws.Range("Q20").Value = 510.49161933274
$ws.Range("R20").Value = 2041.96647733096
$ws.Range("S20").Value = 0.01840997847332952
$ws.Range("T20").Value = 0.01195169536103635
$ws.Range("G21").Value = 11.3729585
$ws.Range("H21").Value = 22.745917
$ws.Range("I21").Value = 0.2058927176688748
$ws.Range("J21").Value = 0.1893032314851988
$ws.Range("O21").Value = 0.06989553179102848
$ws.Range("P21").Value = 0.07402865606362163
$ws.Range("Q21").Value = 399.0485486842342
$ws.Range("R21").Value = 2394.291292105405
$ws.Range("S21").Value = 0.01439098099336609
$ws.Range("T21").Value = 0.01401386381534993
$ws.Range("G22").Value = 11.3729585
$ws.Range("H22").Value = 22.745917
$ws.Range("I22").Value = 0.2058927176688748
$ws.Range("J22").Value = 0.1893032314851988
$ws.Range("M22").Value = 29.15707
$ws.Range("N22").Value = 87.47121
$ws.Range("O22").Value = 0.05808192634815011
$ws.Range("P22").Value = 0.06151647807752574
$ws.Range("Q22").Value = 331.602147091595
$ws.Range("R22").Value = 1989.61288254957
$ws.Range("S22").Value = 0.01195864566326405
$ws.Range("T22").Value = 0.01164526808966401
$ws.Range("G23").Value = 11.3729585
$ws.Range("H23").Value = 22.745917
$ws.Range("I23").Value = 0.2058927176688748
$ws.Range("J23").Value = 0.1893032314851988
$ws.Range("M23").Value = 39.1954995
$ws.Range("N23").Value = 78.39099899999999
$ws.Range("O23").Value = 0.07807883697291786
$ws.Range("P23").Value = 0.05513057577983479
$ws.Range("Q23").Value = 445.7687892002707
$ws.Range("R23").Value = 1783.075156801083
$ws.Range("S23").Value = 0.01607586393677908
$ws.Range("T23").Value = 0.01043639614876236
$ws.Range("G24").Value = 11.3729585
$ws.Range("H24").Value = 22.745917
$ws.Range("I24").Value = 0.2058927176688748
$ws.Range("J24").Value = 0.1893032314851988
$ws.Range("M24").Value = 270.2169853333333
$ws.Range("N24").Value = 810.650956
$ws.Range("O24").Value = 0.5382819000726007
$ws.Range("P24").Value = 0.5701120604516535
$ws.Range("Q24").Value = 3073.166560191108
$ws.Range("R24").Value = 18438.99936114665
$ws.Range("S24").Value = 0.1108283232779135
$ws.Range("T24").Value = 0.107924055352183
$ws.Range("G25").Value = 11.3729585
$ws.Range("H25").Value = 22.745917
$ws.Range("I25").Value = 0.2058927176688748
$ws.Range("J25").Value = 0.1893032314851988
$ws.Range("M25").Value = 83.455535
$ws.Range("N25").Value = 250.366605
$ws.Range("O25").Value = 0.1662464108092982
$ws.Range("P25").Value = 0.1760770402950531
$ws.Range("Q25").Value = 949.1363361502974
$ws.Range("R25").Value = 5694.818016901784
$ws.Range("S25").Value = 0.03422892532422261
$ws.Range("T25").Value = 0.03333195271820311
$ws.Range("G26").Value = 2.577819333333334
$ws.Range("H26").Value = 7.733458000000001
$ws.Range("I26").Value = 0.04666808800888241
$ws.Range("J26").Value = 0.0643618188686375
$ws.Range("M26").Value = 44.88644
$ws.Range("N26").Value = 89.77288
$ws.Range("O26").Value = 0.08941539400600468
$ws.Range("P26").Value = 0.06313518933231116
$ws.Range("Q26").Value = 115.7091328365067
$ws.Range("R26").Value = 694.2547970190401
$ws.Range("S26").Value = 0.004172845476821122
$ws.Range("T26").Value = 0.004063495620043345
$ws.Range("G27").Value = 2.577819333333334
$ws.Range("H27").Value = 7.733458000000001
$ws.Range("I27").Value = 0.04666808800888241
$ws.Range("J27").Value = 0.0643618188686375
$ws.Range("O27").Value = 0.06989553179102848
$ws.Range("P27").Value = 0.07402865606362163
$ws.Range("Q27").Value = 90.44920578377446
$ws.Range("R27").Value = 814.0428520539701
$ws.Range("S27").Value = 0.003261890829051355
$ws.Range("T27").Value = 0.004764618952655479
$ws.Range("G28").Value = 2.577819333333334
$ws.Range("H28").Value = 7.733458000000001
$ws.Range("I28").Value = 0.04666808800888241
$ws.Range("J28").Value = 0.0643618188686375
$ws.Range("M28").Value = 29.15707
$ws.Range("N28").Value = 87.47121
$ws.Range("O28").Value = 0.05808192634815011
$ws.Range("P28").Value = 0.06151647807752574
$ws.Range("Q28").Value = 75.16165874935335
$ws.Range("R28").Value = 676.4549287441801
$ws.Range("S28").Value = 0.002710572450540895
$ws.Range("T28").Value = 0.003959312419462221
$ws.Range("G29").Value = 2.577819333333334
$ws.Range("H29").Value = 7.733458000000001
$ws.Range("I29").Value = 0.04666808800888241
$ws.Range("J29").Value = 0.0643618188686375
$ws.Range("M29").Value = 39.1954995
$ws.Range("N29").Value = 78.39099899999999
$ws.Range("O29").Value = 0.07807883697291786
$ws.Range("P29").Value = 0.05513057577983479
$ws.Range("Q29").Value = 101.038916390757
$ws.Range("R29").Value = 606.233498344542
$ws.Range("S29").Value = 0.003643790035483312
$ws.Range("T29").Value = 0.00354830413246542
$ws.Range("G30").Value = 2.577819333333334
$ws.Range("H30").Value = 7.733458000000001
$ws.Range("I30").Value = 0.04666808800888241
$ws.Range("J30").Value = 0.0643618188686375
$ws.Range("M30").Value = 270.2169853333333
$ws.Range("N30").Value = 810.650956
$ws.Range("O30").Value = 0.5382819000726007
$ws.Range("P30").Value = 0.5701120604516535
$ws.Range("Q30").Value = 696.5705689873164
$ws.Range("R30").Value = 6269.135120885848
$ws.Range("S30").Value = 0.02512058708617657
$ws.Range("T30").Value = 0.03669344916961503
$ws.Range("G31").Value = 2.577819333333334
$ws.Range("H31").Value = 7.733458000000001
$ws.Range("I31").Value = 0.04666808800888241
$ws.Range("J31").Value = 0.0643618188686375
$ws.Range("M31").Value = 83.455535
$ws.Range("N31").Value = 250.366605
$ws.Range("O31").Value = 0.1662464108092982
$ws.Range("P31").Value = 0.1760770402950531
$ws.Range("Q31").Value = 215.1332915966767
$ws.Range("R31").Value = 1936.19962437009
$ws.Range("S31").Value = 0.007758402130809146
$ws.Range("T31").Value = 0.01133263857439599
$ws.Range("G32").Value = 1.820947666666666
$ws.Range("H32").Value = 5.462842999999999
$ws.Range("I32").Value = 0.0329659045025792
$ws.Range("J32").Value = 0.04546459186482997
$ws.Range("M32").Value = 44.88644
$ws.Range("N32").Value = 89.77288
$ws.Range("O32").Value = 0.08941539400600468
$ws.Range("P32").Value = 0.06313518933231116
$ws.Range("Q32").Value = 81.73585818297332
$ws.Range("R32").Value = 490.41514909784
$ws.Range("S32").Value = 0.002947659339862442
$ws.Range("T32").Value = 0.002870415615302294
$ws.Range("G33").Value = 1.820947666666666
$ws.Range("H33").Value = 5.462842999999999
$ws.Range("I33").Value = 0.0329659045025792
$ws.Range("J33").Value = 0.04546459186482997
$ws.Range("O33").Value = 0.06989553179102848
$ws.Range("P33").Value = 0.07402865606362163
$ws.Range("Q33").Value = 63.89248000977721
$ws.Range("R33").Value = 575.032320087995
$ws.Range("S33").Value = 0.002304169426180033
$ws.Range("T33").Value = 0.003365682634234428
$ws.Range("G34").Value = 1.820947666666666
$ws.Range("H34").Value = 5.462842999999999
$ws.Range("I34").Value = 0.0329659045025792
$ws.Range("J34").Value = 0.04546459186482997
$ws.Range("M34").Value = 29.15707
$ws.Range("N34").Value = 87.47121
$ws.Range("O34").Value = 0.05808192634815011
$ws.Range("P34").Value = 0.06151647807752574
$ws.Range("Q34").Value = 53.09349858333666
$ws.Range("R34").Value = 477.84148725003
$ws.Range("S34").Value = 0.001914723237318955
$ws.Range("T34").Value = 0.002796821568756468
$ws.Range("G35").Value = 1.820947666666666
$ws.Range("H35").Value = 5.462842999999999
$ws.Range("I35").Value = 0.0329659045025792
$ws.Range("J35").Value = 0.04546459186482997
$ws.Range("M35").Value = 39.1954995
$ws.Range("N35").Value = 78.39099899999999
$ws.Range("O35").Value = 0.07807883697291786
$ws.Range("P35").Value = 0.05513057577983479
$ws.Range("Q35").Value = 71.37295335835948
$ws.Range("R35").Value = 428.2377201501569
$ws.Range("S35").Value = 0.00257393948332166
$ws.Range("T35").Value = 0.002506489127103269
$ws.Range("G36").Value = 1.820947666666666
$ws.Range("H36").Value = 5.462842999999999
$ws.Range("I36").Value = 0.0329659045025792
$ws.Range("J36").Value = 0.04546459186482997
$ws.Range("M36").Value = 270.2169853333333
$ws.Range("N36").Value = 810.650956
$ws.Range("O36").Value = 0.5382819000726007
$ws.Range("P36").Value = 0.5701120604516535
$ws.Range("Q36").Value = 492.0509889364341
$ws.Range("R36").Value = 4428.458900427907
$ws.Range("S36").Value = 0.01774494971326023
$ws.Range("T36").Value = 0.0259199121456517
$ws.Range("G37").Value = 1.820947666666666
$ws.Range("H37").Value = 5.462842999999999
$ws.Range("I37").Value = 0.0329659045025792
$ws.Range("J37").Value = 0.04546459186482997
$ws.Range("M37").Value = 83.455535
$ws.Range("N37").Value = 250.366605
$ws.Range("O37").Value = 0.1662464108092982
$ws.Range("P37").Value = 0.1760770402950531
$ws.Range("Q37").Value = 151.9681617286683
$ws.Range("R37").Value = 1367.713455558015
$ws.Range("S37").Value = 0.005480463302635873
$ws.Range("T37").Value = 0.008005270773781808
